# CONVERT AND RESIZE IMAGE FOLDERS.docx -- "Updates from testing" edit
#
# Applies a series of targeted Find/Replace operations (scoped to the
# relevant paragraph so that unrelated, repeated words elsewhere in the
# document are not touched) plus one new paragraph (with a hyperlink)
# describing the technology used to build CRIMP.

$d = $word.ActiveDocument

# NOTE: this PowerShell COM host does not honour default parameter
# values for omitted arguments (they come through as blank/false), so
# $matchCase is always passed explicitly below.
function Replace-InRange($range, $old, $new, $matchCase) {
    $range.Find.ClearFormatting()
    $range.Find.Execute($old, $matchCase, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------
# 1. "Same as one 1; you want to..." -> "You want to ... those photographs"
# ---------------------------------------------------------------------
$p = $d.Paragraphs(8).Range
Replace-InRange $p `
    "Same as one 1; you want to add the photograph file name to the photograph so that someone wanting a full-size copy can tell you which photograph(s) they want." `
    "You want to add the photograph file name to those photographs so that someone wanting a full-size copy can tell you which photograph(s) they want." `
    $true

# ---------------------------------------------------------------------
# 2. History paragraph: typo fixes + re-worded "extends the capabilities" sentence
# ---------------------------------------------------------------------
$p = $d.Paragraphs(14).Range
Replace-InRange $p "could covert their images" "could convert their images" $true

$p = $d.Paragraphs(14).Range
Replace-InRange $p "such as PhotoShop or Paint.Net" "such as Photoshop or Paint.Net" $true

$p = $d.Paragraphs(14).Range
Replace-InRange $p `
    "CRIMP extends the capabilities of its predecessor to convert files between the most common formats used by photographers. It can convert files to and from BMP, GIF, TIFF, JPG. and WEBP and resize them at the same time." `
    "CRIMP extends the capabilities of its predecessor to the conversion of files from and to the most common formats used by photographers. It can convert BMP, GIF, TIFF, JPG. and WEBP files and resize them at the same time." `
    $true

# ---------------------------------------------------------------------
# 3. "Output Folder " button label gains padding spaces
# ---------------------------------------------------------------------
$p = $d.Paragraphs(31).Range
Replace-InRange $p "Output Folder " "  Output Folder  " $true

# ---------------------------------------------------------------------
# 4. "Select Folders" button label -> "Select" (padded), "Folders" dropped
# ---------------------------------------------------------------------
$p = $d.Paragraphs(35).Range
Replace-InRange $p "Select Folders" " Select " $true

# ---------------------------------------------------------------------
# 5. "Check All" button label gains padding spaces
# ---------------------------------------------------------------------
$p = $d.Paragraphs(36).Range
Replace-InRange $p "Check All" " Check All " $true

# ---------------------------------------------------------------------
# 6. Drop the "This lets you click the Process Folders button" sentence
# ---------------------------------------------------------------------
$p = $d.Paragraphs(40).Range
Replace-InRange $p `
    " button. This lets you click the Process Folders button. CRIMP will check" `
    " button. CRIMP will check" `
    $true

# ---------------------------------------------------------------------
# 7. "clicking the Show Stats button" -> "switching to the Show Stats tab"
# ---------------------------------------------------------------------
$p = $d.Paragraphs(41).Range
Replace-InRange $p "run by clicking the " "run by switching to the " $true

$p = $d.Paragraphs(41).Range
Replace-InRange $p " button." " tab." $true

# ---------------------------------------------------------------------
# 8. "By clicking the Error Log button" -> "By switching to the Error Log tab"
# ---------------------------------------------------------------------
$p = $d.Paragraphs(42).Range
Replace-InRange $p "run. By clicking the " "run. By switching to the " $true

$p = $d.Paragraphs(42).Range
Replace-InRange $p " button. Typical" " tab. Typical" $true

# ---------------------------------------------------------------------
# 9. New "Technical Information" paragraph with a Magick.Net hyperlink,
#    inserted right after the "Technical Information" heading.
# ---------------------------------------------------------------------
$headingIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Technical Information") {
        $headingIndex = $i
        break
    }
}

$heading = $d.Paragraphs($headingIndex).Range
$heading.InsertParagraphAfter()

$newPara = $d.Paragraphs($headingIndex + 1)
$newPara.Style = "HelpText"
$newPara.Range.Text = "CRIMP is a WPF application developed using Microsoft Visual Studio 2022. It uses Magick.Net to read, convert, resize and save individual files. The source code is available on GitHub."

$linkRange = $d.Paragraphs($headingIndex + 1).Range
$linkRange.Find.ClearFormatting()
$linkRange.Find.Execute("Magick.Net", $true) | Out-Null
$d.Hyperlinks.Add($linkRange, "https://github.com/dlemstra/Magick.NET", $null, $null, "Magick.Net") | Out-Null

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
